$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$rowData = @{
    2 = @{ D="300.93"; E="-0.59%"; G="3" }
    3 = @{ D="38.22"; E="9.02%"; G="3" }
    4 = @{ E="-3.30%"; G="3" }
    5 = @{ D="0.07708"; E="-0.82%"; G="3" }
    6 = @{ D="2.191"; E="-6.08%"; G="3" }
    7 = @{ D="7.969"; E="-0.80%"; G="3" }
    8 = @{ D="3.994"; E="1.10%"; G="3" }
    9 = @{ D="0.9167"; E="-1.66%"; G="3" }
    10 = @{ D="0.08998"; E="-9.28%"; G="3" }
    11 = @{ D="0.1786"; E="-0.40%"; G="3" }
    12 = @{ D="0.08462"; E="-1.87%"; G="3" }
    13 = @{ D="0.03536"; E="6.46%"; G="3" }
    14 = @{ D="0.09936"; E="0.14%"; G="3" }
    15 = @{ D="0.001471"; E="-1.96%"; G="3" }
    16 = @{ D="0.005658"; E="-2.85%"; G="3" }
    17 = @{ D="3.479"; E="0.48%"; G="3" }
    18 = @{ D="2.223"; E="4.15%"; G="3" }
    19 = @{ E="2.86%"; G="3" }
    20 = @{ D="0.1301"; E="-2.47%"; G="3" }
    21 = @{ D="4.559"; E="5.63%"; G="3" }
    22 = @{ E="-2.69%"; G="3" }
    23 = @{ D="0.04661"; E="1.01%"; G="3" }
    24 = @{ D="0.001231"; E="1.17%"; G="3" }
    25 = @{ D="0.004433"; E="1.28%"; G="3" }
    26 = @{ D="0.0001303"; E="0.08%"; G="3" }
    27 = @{ D="0.0004758"; E="40.15%"; G="3" }
    28 = @{ G="3" }
    29 = @{ G="3" }
    30 = @{ G="3" }
    31 = @{ G="3" }
    32 = @{ G="3" }
    33 = @{ G="3" }
    34 = @{ G="3" }
    35 = @{ G="3" }
    36 = @{ G="3" }
    37 = @{ G="3" }
    38 = @{ G="3" }
    39 = @{ D="0.01739"; E="-2.67%"; G="3" }
    40 = @{ D="0.04670"; E="-2.78%"; G="3" }
    41 = @{ D="0.007855"; E="0.81%"; G="3" }
    42 = @{ D="0.1386"; E="-1.79%"; G="3" }
    43 = @{ D="0.007683"; E="2.56%"; G="3" }
    44 = @{ D="0.002295"; E="9.51%"; G="3" }
    45 = @{ D="0.009557"; E="1.07%"; G="3" }
    46 = @{ D="0.00006059"; E="-0.93%"; G="3" }
    47 = @{ E="0.07%"; G="3" }
    48 = @{ D="8.520"; E="190.26%"; G="3" }
    49 = @{ E="34.73%"; G="3" }
    50 = @{ D="0.00002105"; E="0.07%"; G="3" }
    51 = @{ D="0.0002004"; E="0.07%"; G="3" }
}

foreach ($r in $rowData.Keys) {
    $cellsForRow = $rowData[$r]
    foreach ($c in $cellsForRow.Keys) {
        $addr = "$c$r"
        Set-TextValue $ws.Range($addr) $cellsForRow[$c]
    }
}

Write-Host "Applied $($rowData.Count) row updates"